# Update Name of Algo
# Apply corrected values to the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.41060000000003
$ws.Range("C5").Value = -14.14190000000001
$ws.Range("C9").Value = -11.95440000000001
$ws.Range("C11").Value = -13.53450000000001
$ws.Range("A21").Value = -20.83070000000001
$ws.Range("C21").Value = -11.3248
$ws.Range("A23").Value = -21.49450000000003
$ws.Range("A25").Value = -22.39150000000004

$wb.Save()
